# Auto-generated edit script applying the Ultros_Profits diff
# Updates H/I/J/K/L/M/N value cells across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2049.158
$ws.Range("I15").Value = 2049.158
$ws.Range("K15").Value = 6147.474
$ws.Range("M15").Value = -5978.474
$ws.Range("H74").Value = 9300.333000000001
$ws.Range("I74").Value = 7962.875
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 7962.875
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -7026.875
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 9300.333000000001
$ws.Range("I77").Value = 7962.875
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 39814.375
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -35134.375
$ws.Range("N77").Value = -109360
$ws.Range("H80").Value = 3351.96
$ws.Range("I80").Value = 982.375
$ws.Range("K80").Value = 2947.125
$ws.Range("M80").Value = -1949.125
$ws.Range("H83").Value = 3351.96
$ws.Range("I83").Value = 982.375
$ws.Range("K83").Value = 8841.375
$ws.Range("M83").Value = -3849.375
$ws.Range("H92").Value = 974.1429000000001
$ws.Range("I92").Value = 835.3889
$ws.Range("K92").Value = 835.3889
$ws.Range("M92").Value = 412.6111
$ws.Range("H104").Value = 181.5
$ws.Range("I104").Value = 181.5
$ws.Range("K104").Value = 544.5
$ws.Range("M104").Value = 1202.5
$ws.Range("H108").Value = 32675
$ws.Range("J108").Value = 32675
$ws.Range("L108").Value = 32675
$ws.Range("N108").Value = -40355
$ws.Range("H113").Value = 9937.179
$ws.Range("I113").Value = 6605.25
$ws.Range("K113").Value = 6605.25
$ws.Range("M113").Value = -3351.25
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140
$ws.Range("H138").Value = 3392.6
$ws.Range("I138").Value = 1648.5
$ws.Range("K138").Value = 4945.5
$ws.Range("M138").Value = 194.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 555000000
$ws.Range("I5").Value = 110000000
$ws.Range("K5").Value = 110000000
$ws.Range("M5").Value = -109999888
$ws.Range("H45").Value = 3432.625
$ws.Range("I45").Value = 3142.3
$ws.Range("J45").Value = 3916.5
$ws.Range("K45").Value = 3142.3
$ws.Range("L45").Value = 3916.5
$ws.Range("M45").Value = -2765.3
$ws.Range("N45").Value = -4670.5
$ws.Range("H61").Value = 4702.0835
$ws.Range("I61").Value = 4603.1113
$ws.Range("K61").Value = 4603.1113
$ws.Range("M61").Value = -4391.1113
$ws.Range("H63").Value = 7240.3
$ws.Range("I63").Value = 4343.4287
$ws.Range("J63").Value = 13999.667
$ws.Range("K63").Value = 4343.4287
$ws.Range("L63").Value = 13999.667
$ws.Range("M63").Value = -3657.4287
$ws.Range("N63").Value = -15371.667
$ws.Range("H66").Value = 7240.3
$ws.Range("I66").Value = 4343.4287
$ws.Range("J66").Value = 13999.667
$ws.Range("K66").Value = 21717.1435
$ws.Range("L66").Value = 69998.33499999999
$ws.Range("M66").Value = -18285.1435
$ws.Range("N66").Value = -76862.33499999999
$ws.Range("H97").Value = 2436.647
$ws.Range("I97").Value = 1130
$ws.Range("K97").Value = 1130
$ws.Range("M97").Value = -634
$ws.Range("H102").Value = 4371.7617
$ws.Range("I102").Value = 2369.1875
$ws.Range("K102").Value = 2369.1875
$ws.Range("M102").Value = -747.1875
$ws.Range("H122").Value = 4698.7856
$ws.Range("I122").Value = 4217.476
$ws.Range("K122").Value = 12652.428
$ws.Range("M122").Value = -10202.428
$ws.Range("H132").Value = 31252494
$ws.Range("I132").Value = 41668940
$ws.Range("K132").Value = 125006820
$ws.Range("M132").Value = -125004290
$ws.Range("H136").Value = 4702.0835
$ws.Range("I136").Value = 4603.1113
$ws.Range("K136").Value = 13809.3339
$ws.Range("M136").Value = -11259.3339
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 555000000
$ws.Range("I4").Value = 110000000
$ws.Range("K4").Value = 110000000
$ws.Range("M4").Value = -109999885
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H99").Value = 40530.617
$ws.Range("I99").Value = 1132
$ws.Range("K99").Value = 1132
$ws.Range("M99").Value = 366
$ws.Range("H134").Value = 3483
$ws.Range("I134").Value = 2779.6
$ws.Range("K134").Value = 8338.799999999999
$ws.Range("M134").Value = -5803.799999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 520.9286
$ws.Range("I22").Value = 599.3333
$ws.Range("K22").Value = 599.3333
$ws.Range("M22").Value = -249.3333
$ws.Range("H86").Value = 39626
$ws.Range("I86").Value = 55757.168
$ws.Range("K86").Value = 55757.168
$ws.Range("M86").Value = -54634.168
$ws.Range("H89").Value = 39626
$ws.Range("I89").Value = 55757.168
$ws.Range("K89").Value = 278785.84
$ws.Range("M89").Value = -273169.84
$ws.Range("H105").Value = 36111652
$ws.Range("I105").Value = 18518938
$ws.Range("K105").Value = 18518938
$ws.Range("M105").Value = -18517191
$ws.Range("H132").Value = 3586.4348
$ws.Range("I132").Value = 2471.889
$ws.Range("K132").Value = 7415.667
$ws.Range("M132").Value = -4885.667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 496.83334
$ws.Range("I5").Value = 446.2
$ws.Range("K5").Value = 1338.6
$ws.Range("M5").Value = -1226.6
$ws.Range("H131").Value = 2207.2
$ws.Range("I131").Value = 1657.8572
$ws.Range("J131").Value = 2687.875
$ws.Range("K131").Value = 4973.571599999999
$ws.Range("L131").Value = 8063.625
$ws.Range("M131").Value = 66.42840000000069
$ws.Range("N131").Value = -18143.625
$ws.Range("H135").Value = 496.83334
$ws.Range("I135").Value = 446.2
$ws.Range("K135").Value = 4015.8
$ws.Range("M135").Value = -1480.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3929.8572
$ws.Range("I97").Value = 387.22726
$ws.Range("K97").Value = 387.22726
$ws.Range("M97").Value = 108.77274
$ws.Range("H113").Value = 7316.273
$ws.Range("I113").Value = 4497.6665
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 4497.6665
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = -2327.6665
$ws.Range("N113").Value = -24340
$ws.Range("H126").Value = 2271.4285
$ws.Range("I126").Value = 1881.8182
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 5645.4546
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -3175.4546
$ws.Range("N126").Value = -16040
$ws.Range("H132").Value = 9006.346
$ws.Range("I132").Value = 8041.7
$ws.Range("K132").Value = 24125.1
$ws.Range("M132").Value = -21595.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2879.3635
$ws.Range("I46").Value = 1433.3334
$ws.Range("J46").Value = 3421.625
$ws.Range("K46").Value = 1433.3334
$ws.Range("L46").Value = 3421.625
$ws.Range("M46").Value = -1245.3334
$ws.Range("N46").Value = -3797.625
$ws.Range("H93").Value = 1890
$ws.Range("I93").Value = 1744.4445
$ws.Range("K93").Value = 1744.4445
$ws.Range("M93").Value = -496.4445000000001
$ws.Range("H100").Value = 128900
$ws.Range("I100").Value = 282025
$ws.Range("J100").Value = 6400
$ws.Range("K100").Value = 282025
$ws.Range("L100").Value = 6400
$ws.Range("M100").Value = -281484
$ws.Range("N100").Value = -7482
$ws.Range("H139").Value = 99992
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3182.2
$ws.Range("I96").Value = 3477.75
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 3477.75
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -2104.75
$ws.Range("N96").Value = -4746
$ws.Range("H100").Value = 726.1111
$ws.Range("I100").Value = 788.8333
$ws.Range("J100").Value = 600.6667
$ws.Range("K100").Value = 1577.6666
$ws.Range("L100").Value = 1201.3334
$ws.Range("M100").Value = -1036.6666
$ws.Range("N100").Value = -2283.3334
$ws.Range("H119").Value = 68142.5
$ws.Range("J119").Value = 68142.5
$ws.Range("L119").Value = 68142.5
$ws.Range("N119").Value = -77818.5
$ws.Range("H139").Value = 74996

